$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add new column F with header and values
$ws.Range("F1").Value = "25_03_2024"
$ws.Range("F2").Value = 1190
$ws.Range("F3").Value = 1122
$ws.Range("F4").Value = 1602
$ws.Range("F5").Value = 220

# Update selection to match target (F6)
$ws.Range("F6").Select()
